$d = $word.ActiveDocument
$d.Content.Find.Execute("sản phẩm", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Book", 2)
